# Generate Report for Handback
# Adds a new handback entry (d6199722-e2fb-4701-ad20-c13c4e29cf31.md) as row 4
# to the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$newFile   = "d6199722-e2fb-4701-ad20-c13c4e29cf31.md"
$newPath   = "e2e\d6199722-e2fb-4701-ad20-c13c4e29cf31.md"
$newStatus = "Handed back: in sync with en-US"
$ext       = ".md"

$zhXlf     = "d6199722-e2fb-4701-ad20-c13c4e29cf31.327012bcaf7f839c9096b3b08d52309cd7db483d.zh-cn.xlf"
$deXlf     = "d6199722-e2fb-4701-ad20-c13c4e29cf31.327012bcaf7f839c9096b3b08d52309cd7db483d.de-de.xlf"

$zhHandoffDate   = "'2016-09-07 07:53:11"
$zhHandbackDate  = "'2016-09-07 07:54:09"
$deHandoffDate   = "'2016-09-07 07:53:24"
$deHandbackDate  = "'2016-09-07 07:54:30"
$overviewDate    = "'2016-09-07 07:53:24"

# ---------------------------------------------------------------------------
# Sheet "Overview" - row 4
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(4, 1).Value = $newFile
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newPath)
$ov.Cells.Item(4, 3).Value = $ext
$ov.Cells.Item(4, 5).Value = $newStatus
$ov.Cells.Item(4, 6).Value = $newStatus
$ov.Cells.Item(4, 7).Value = $overviewDate
$ov.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" - row 4
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile)
$zh.Cells.Item(4, 2).Value = $ext
$zh.Cells.Item(4, 3).Value = $newStatus
$zh.Cells.Item(4, 4).Value = "e2e"
$zh.Cells.Item(4, 5).Value = "ht"
$zh.Cells.Item(4, 6).Value = "'True"
$zh.Cells.Item(4, 7).Value = $zhXlf
$zh.Cells.Item(4, 8).Value = $zhHandoffDate
$zh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Hyperlinks.Add($zh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile)
$zh.Cells.Item(4, 10).Value = $zhXlf
$zh.Cells.Item(4, 11).Value = $zhHandbackDate
$zh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item(4, 12).Value = "'"
$zh.Cells.Item(4, 13).Value = "'True"
$zh.Cells.Item(4, 14).Value = "'"
$zh.Cells.Item(4, 15).Value = "'False"
$zh.Cells.Item(4, 16).Value = "'"

# ---------------------------------------------------------------------------
# Sheet "de-de" - row 4
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile)
$de.Cells.Item(4, 2).Value = $ext
$de.Cells.Item(4, 3).Value = $newStatus
$de.Cells.Item(4, 4).Value = "e2e"
$de.Cells.Item(4, 5).Value = "ht"
$de.Cells.Item(4, 6).Value = "'True"
$de.Cells.Item(4, 7).Value = $deXlf
$de.Cells.Item(4, 8).Value = $deHandoffDate
$de.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Hyperlinks.Add($de.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0000000000000000000000000000000000000000/e2e/$newFile", "", "", $newFile)
$de.Cells.Item(4, 10).Value = $deXlf
$de.Cells.Item(4, 11).Value = $deHandbackDate
$de.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item(4, 12).Value = "'"
$de.Cells.Item(4, 13).Value = "'True"
$de.Cells.Item(4, 14).Value = "'"
$de.Cells.Item(4, 15).Value = "'False"
$de.Cells.Item(4, 16).Value = "'"

# ---------------------------------------------------------------------------
# Extend the tables + dimensions to include the new row
# ---------------------------------------------------------------------------
$ov.ListObjects.Item(1).Resize($ov.Range("A1:G4"))
$zh.ListObjects.Item(1).Resize($zh.Range("A1:P4"))
$de.ListObjects.Item(1).Resize($de.Range("A1:P4"))
